# edit.ps1
#
# Applies two changes described by the commit "added image for Number
# Theory (part 2) and edited Eulers_Theorem.json":
#
#   1. The cached text of the `datetimeFigureOut` field ("today" date
#      placeholder) is refreshed from 6/18/2024 to 3/7/2025 on the
#      slide master and on every one of the 11 slide layouts.
#   2. On slide 2 ("NUMBER" / "THEORY" title card), the second line of
#      text is updated from "THEORY" to "THEORY 2".
#
# Helper: overwrite a shape's whole text via a full-length Characters()
# sub-range rather than TextRange.Text directly. For ordinary runs this
# produces an identical result to setting .Text, but it edits the
# existing run(s) in place (matching how PowerPoint itself keeps a
# single run when only its characters change) instead of collapsing
# every paragraph in the shape into one run.
function Set-ShapeText {
    param($Shape, [string]$NewText)
    $tr = $Shape.TextFrame.TextRange
    $len = $tr.Length
    $sub = $tr.Characters(1, $len)
    $sub.Text = $NewText
}

$p = $ppt.ActivePresentation

# --- 1. Date placeholder: 6/18/2024 -> 3/7/2025 -------------------------

$newDate = "3/7/2025"

# Slide master's own "Date Placeholder" shape.
Set-ShapeText $p.SlideMaster.Shapes.Item(3) $newDate

# Every slide layout ("custom layout") also carries its own cached copy
# of the date placeholder text.
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $shapes = $layout.Shapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shape = $shapes.Item($j)
        if ($shape.Name -like "Date Placeholder*") {
            Set-ShapeText $shape $newDate
        }
    }
}

# --- 2. Slide 2 title text: "THEORY" -> "THEORY 2" ----------------------

$slide2 = $p.Slides.Item(2)
$titleShape = $slide2.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange

$idx = $tr.Text.IndexOf("THEORY")
if ($idx -ge 0) {
    $word = $tr.Characters($idx + 1, 6)
    $word.Text = "THEORY 2"
}
